# semester-config.py creates new semester entries in the DB - this adds the
# matching "semester" worksheet (lookup table of YEAR / SEASON / COURSE_NO)
# to the workbook, after the existing "students" sheet.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end (after the last current sheet) and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "semester"

# Header row
$ws.Range("A1").Value = "YEAR"
$ws.Range("B1").Value = "SEASON"
$ws.Range("C1").Value = "COURSE_NO"

# Data rows
$ws.Range("A2").Value = 2016
$ws.Range("B2").Value = "Fall"
$ws.Range("C2").Value = "P532"

$ws.Range("A3").Value = 2016
$ws.Range("B3").Value = "Spring"
$ws.Range("C3").Value = "P632"

$ws.Range("A4").Value = 2017
$ws.Range("B4").Value = "Spring"
$ws.Range("C4").Value = "P532"

# Auto-size the COURSE_NO column, matching the "bestFit" columns already used
# on the "students" sheet.
$ws.Columns.Item(3).AutoFit()

# The previous selection/cursor on "students" moves off of it (it is no
# longer the active tab) - restore its last-known selection.
$students = $wb.Worksheets.Item("students")
$null = $students.Range("C39").Select()

# "semester" becomes the active (selected) sheet/tab with its own selection.
$null = $ws.Activate()
$null = $ws.Range("C4").Select()
